$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").Value = 13.102169094515716
$ws.Range("C2").Value = -0.49387760034906591
$ws.Range("D2").Value = 0.50889455340507084
$ws.Range("E2").Value = 0.24106549700383084

$ws.Range("B3").Value = 4.0443387243297195
$ws.Range("C3").Value = 8.8176957739006525
$ws.Range("D3").Value = 3.0261117330107936
$ws.Range("E3").Value = -1.0589812866004138

[void]$ws.Range("B1:E3").Select()
